$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell E8 content ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection change recorded in the saved file
$ws.Range("E8").Select()
